# Q1, partB3,4 hesaplamalar dışında tamamlandı
# Update the last few data rows (B-H curve of core), drop the stray extra
# row, add the B(mT)-H(At/m) line chart and move the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# --- finish filling in the H (col C) values for the last data points and
#     correct the trailing A-column (u-derived B) values ---
$ws.Range("A13").Value = 520
$ws.Range("C13").Value = 90

$ws.Range("A14").Value = 530
$ws.Range("C14").Value = 100

$ws.Range("A15").Value = 530
$ws.Range("C15").Value = 1050

# row 16 was a leftover duplicate - clear it so the sheet shrinks back to A1:C15
$ws.Range("A16").Value = ""

# --- add the "B(mT)-H(At/m) curve" line chart, plotting B (A2:A15) against
#     H (C2:C15) ---
$co = $ws.ChartObjects().Add(123825, 28575, 5238750, 3552824)
$chart = $co.Chart
$chart.ChartType = 63

$chart.SetSourceData($ws.Range("A2:A15")) | Out-Null

$series = $chart.SeriesCollection(1)
$series.Name = "=Sayfa1!`$A`$1"
$series.XValues = $ws.Range("C2:C15")
$series.Values = $ws.Range("A2:A15")
$series.MarkerStyle = -4142

$chart.HasTitle = $true
$chart.ChartTitle.Text = "B(mT)-H(At/m) curve"

$chart.HasLegend = $true
$chart.Legend.Position = -4152

# --- move the active selection to D11, matching the saved workbook state ---
$ws.Range("D11").Select() | Out-Null
